$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns D and E ---
$ws.Range("D1").Value = "Portfel"
$ws.Range("E1").Value = "Koronowirus"

# --- Update amounts (column D) and infection status (column E) for existing rows ---
$ws.Range("D2").Value = 50.0
$ws.Range("E2").Value = "brak"

$ws.Range("D3").Value = 250.0
$ws.Range("E3").Value = "brak"

$ws.Range("D4").Value = 2100.0
$ws.Range("E4").Value = "negatywny"

$ws.Range("D5").Value = 50.0
$ws.Range("E5").Value = "brak"

$ws.Range("D6").Value = 150.0
$ws.Range("E6").Value = "brak"

$ws.Range("D7").Value = 2100.0
$ws.Range("E7").Value = "pozytywny"

# --- Insert a new row 9 with the old row-8 (Justyna) data, shifted down and updated ---
$ws.Range("A9").Value = "Justyna"
$ws.Range("B9").Value = "Justynowska"
$ws.Range("C9").Value = "'465633636"
$ws.Range("D9").Value = 250.0
$ws.Range("E9").Value = "brak"

# --- Row 8 becomes a duplicate of row 7 (Ela Lewacka, positive result) ---
$ws.Range("A8").Value = "Ela"
$ws.Range("B8").Value = "Lewacka"
$ws.Range("C8").Value = "'53400564"
$ws.Range("D8").Value = 2100.0
$ws.Range("E8").Value = "pozytywny"
